# Apply updated crypto market data (price + 1h volume change) per the
# Fri May  5 07:40:57 UTC 2023 GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.106.28"
Set-TextValue "E2" "  -0.26%  "

# Row 3
Set-TextValue "D3" "1.897.27"
Set-TextValue "E3" "  -0.45%  "

# Row 4
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.30%  "

# Row 5
Set-TextValue "D5" "325.18"
Set-TextValue "E5" "  -0.72%  "

# Row 6
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.34%  "

# Row 7
Set-TextValue "D7" "0.4623"
Set-TextValue "E7" "  -0.35%  "

# Row 8
Set-TextValue "D8" "0.3896"
Set-TextValue "E8" "  -1.40%  "

# Row 9
Set-TextValue "D9" "0.07882"
Set-TextValue "E9" "  -1.03%  "

# Row 10
Set-TextValue "D10" "0.9891"
Set-TextValue "E10" "  -1.06%  "

# Row 11
Set-TextValue "D11" "21.79"
Set-TextValue "E11" "  -2.01%  "

# Row 12
Set-TextValue "D12" "1.881.89"
Set-TextValue "E12" "  -1.40%  "

# Row 13
Set-TextValue "E13" "  -0.64%  "

# Row 14
Set-TextValue "D14" "5.757"
Set-TextValue "E14" "  -0.15%  "

# Row 15
Set-TextValue "D15" "0.07003"
Set-TextValue "E15" "  +0.85%  "

# Row 16
Set-TextValue "D16" "87.99"
Set-TextValue "E16" "  -0.77%  "

# Row 17
Set-TextValue "D17" "1.003"
Set-TextValue "E17" "  -0.19%  "

# Row 18
Set-TextValue "D18" "0.000009980"
Set-TextValue "E18" "  -0.92%  "

# Row 19
Set-TextValue "D19" "17.06"
Set-TextValue "E19" "  -0.63%  "

# Row 20
Set-TextValue "D20" "1.000"
Set-TextValue "E20" "  -0.33%  "

# Row 21
Set-TextValue "D21" "29.125.40"
Set-TextValue "E21" "  -0.29%  "

# Row 22
Set-TextValue "D22" "5.313"
Set-TextValue "E22" "  -0.80%  "

# Row 23
Set-TextValue "D23" "11.10"
Set-TextValue "E23" "  +0.08%  "

# Row 24
Set-TextValue "D24" "2.110.53"
Set-TextValue "E24" "  -1.58%  "

# Row 25
Set-TextValue "D25" "2.107"
Set-TextValue "E25" "  +2.76%  "

# Row 26
Set-TextValue "D26" "155.85"
Set-TextValue "E26" "  -0.53%  "

# Row 27
Set-TextValue "D27" "19.37"
Set-TextValue "E27" "  -0.79%  "

# Row 28
Set-TextValue "D28" "5.946"
Set-TextValue "E28" "  +0.21%  "

# Row 29
Set-TextValue "D29" "118.43"
Set-TextValue "E29" "  -0.46%  "

# Row 30
Set-TextValue "D30" "1.882"
Set-TextValue "E30" "  -5.56%  "

# Row 31
Set-TextValue "D31" "0.09324"
Set-TextValue "E31" "  -0.82%  "

# Row 32
Set-TextValue "D32" "0.8997"
Set-TextValue "E32" "  -2.40%  "

# Row 33
Set-TextValue "D33" "5.248"
Set-TextValue "E33" "  -1.90%  "

# Row 34
Set-TextValue "D34" "1.324"
Set-TextValue "E34" "  -1.78%  "

# Row 35
Set-TextValue "D35" "3.158"
Set-TextValue "E35" "  -3.22%  "

# Row 36
Set-TextValue "B36" "Hedera"
Set-TextValue "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D36" "0.05780"
Set-TextValue "E36" "  -0.74%  "

# Row 37
Set-TextValue "B37" "TrustWalletToken"
Set-TextValue "C37" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D37" "1.176"
Set-TextValue "E37" "  -0.11%  "

# Row 38
Set-TextValue "E38" "  -0.96%  "

# Row 39
Set-TextValue "D39" "1.000"
Set-TextValue "E39" "  -0.27%  "

# Row 40
Set-TextValue "D40" "7.726"
Set-TextValue "E40" "  -3.23%  "

# Row 41
Set-TextValue "D41" "0.5690"
Set-TextValue "E41" "  -1.05%  "

# Row 42
Set-TextValue "D42" "0.1788"
Set-TextValue "E42" "  -0.87%  "

# Row 43
Set-TextValue "D43" "9.698"
Set-TextValue "E43" "  -2.75%  "

# Row 44
Set-TextValue "D44" "11.89"
Set-TextValue "E44" "  -1.11%  "

# Row 45
Set-TextValue "D45" "0.5343"
Set-TextValue "E45" "  -1.59%  "

# Row 46
Set-TextValue "D46" "2.168"
Set-TextValue "E46" "  -2.06%  "

# Row 47
Set-TextValue "D47" "0.07009"
Set-TextValue "E47" "  -1.31%  "

# Row 48
Set-TextValue "D48" "1.853"
Set-TextValue "E48" "  -1.23%  "

# Row 49
Set-TextValue "D49" "2.555"
Set-TextValue "E49" "  -0.28%  "

# Row 50
Set-TextValue "D50" "113.28"
Set-TextValue "E50" "  +1.09%  "

# Row 51
Set-TextValue "E51" "  -0.93%  "
